$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1591.7142
$ws_ALC.Range("J17").Value = 1591.7142
$ws_ALC.Range("L17").Value = 4775.142599999999
$ws_ALC.Range("N17").Value = -5111.142599999999
$ws_ALC.Range("H112").Value = 2361.2273
$ws_ALC.Range("J112").Value = 2430.8096
$ws_ALC.Range("L112").Value = 7292.4288
$ws_ALC.Range("N112").Value = -9508.4288
$ws_ALC.Range("H129").Value = 1108.2354
$ws_ALC.Range("J129").Value = 1169.4333
$ws_ALC.Range("L129").Value = 3508.2999
$ws_ALC.Range("N129").Value = -13508.2999
$ws_ALC.Range("H137").Value = 1421.4468
$ws_ALC.Range("I137").Value = 1378.8889
$ws_ALC.Range("J137").Value = 1478.9
$ws_ALC.Range("K137").Value = 4136.6667
$ws_ALC.Range("L137").Value = 4436.700000000001
$ws_ALC.Range("M137").Value = -1586.6667
$ws_ALC.Range("N137").Value = -9536.700000000001
$ws_ALC.Range("H138").Value = 1912.03
$ws_ALC.Range("I138").Value = 1065.7273
$ws_ALC.Range("J138").Value = 2150.7307
$ws_ALC.Range("K138").Value = 3197.1819
$ws_ALC.Range("L138").Value = 6452.1921
$ws_ALC.Range("M138").Value = 1942.8181
$ws_ALC.Range("N138").Value = -16732.1921
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H13").Value = 171699970
$ws_ARM.Range("J13").Value = 49933.332
$ws_ARM.Range("L13").Value = 49933.332
$ws_ARM.Range("N13").Value = -50221.332
$ws_ARM.Range("H32").Value = 19245.72
$ws_ARM.Range("I32").Value = 18513.477
$ws_ARM.Range("K32").Value = 18513.477
$ws_ARM.Range("M32").Value = -18226.477
$ws_ARM.Range("H61").Value = 1040.5333
$ws_ARM.Range("I61").Value = 1028.1428
$ws_ARM.Range("J61").Value = 1214
$ws_ARM.Range("K61").Value = 1028.1428
$ws_ARM.Range("L61").Value = 1214
$ws_ARM.Range("M61").Value = -816.1428000000001
$ws_ARM.Range("N61").Value = -1638
$ws_ARM.Range("H64").Value = 35000
$ws_ARM.Range("J64").Value = 35000
$ws_ARM.Range("L64").Value = 35000
$ws_ARM.Range("N64").Value = -35496
$ws_ARM.Range("H67").Value = 35000
$ws_ARM.Range("J67").Value = 35000
$ws_ARM.Range("L67").Value = 35000
$ws_ARM.Range("N67").Value = -36716
$ws_ARM.Range("H122").Value = 2741.2
$ws_ARM.Range("I122").Value = 2434
$ws_ARM.Range("J122").Value = 3970
$ws_ARM.Range("K122").Value = 7302
$ws_ARM.Range("L122").Value = 11910
$ws_ARM.Range("M122").Value = -4852
$ws_ARM.Range("N122").Value = -16810
$ws_ARM.Range("H136").Value = 1040.5333
$ws_ARM.Range("I136").Value = 1028.1428
$ws_ARM.Range("J136").Value = 1214
$ws_ARM.Range("K136").Value = 3084.4284
$ws_ARM.Range("L136").Value = 3642
$ws_ARM.Range("M136").Value = -534.4284000000002
$ws_ARM.Range("N136").Value = -8742
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H62").Value = 30000
$ws_BSM.Range("J62").Value = 30000
$ws_BSM.Range("L62").Value = 30000
$ws_BSM.Range("N62").Value = -31372
$ws_BSM.Range("H65").Value = 30000
$ws_BSM.Range("J65").Value = 30000
$ws_BSM.Range("L65").Value = 90000
$ws_BSM.Range("N65").Value = -96864
$ws_BSM.Range("H109").Value = 27251.6
$ws_BSM.Range("J109").Value = 27251.6
$ws_BSM.Range("L109").Value = 27251.6
$ws_BSM.Range("N109").Value = -30025.6
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 33336728
$ws_CRP.Range("I31").Value = 47621604
$ws_CRP.Range("J31").Value = 5355.5557
$ws_CRP.Range("K31").Value = 47621604
$ws_CRP.Range("L31").Value = 5355.5557
$ws_CRP.Range("M31").Value = -47621309
$ws_CRP.Range("N31").Value = -5945.5557
$ws_CRP.Range("H34").Value = 33336728
$ws_CRP.Range("I34").Value = 47621604
$ws_CRP.Range("J34").Value = 5355.5557
$ws_CRP.Range("K34").Value = 47621604
$ws_CRP.Range("L34").Value = 5355.5557
$ws_CRP.Range("M34").Value = -47621402
$ws_CRP.Range("N34").Value = -5759.5557
$ws_CRP.Range("H58").Value = 1819.48
$ws_CRP.Range("I58").Value = 1702.9546
$ws_CRP.Range("K58").Value = 1702.9546
$ws_CRP.Range("M58").Value = -1499.9546
$ws_CRP.Range("H68").Value = 28500
$ws_CRP.Range("I68").Value = 25000
$ws_CRP.Range("K68").Value = 25000
$ws_CRP.Range("M68").Value = -24251
$ws_CRP.Range("H71").Value = 28500
$ws_CRP.Range("I71").Value = 25000
$ws_CRP.Range("K71").Value = 75000
$ws_CRP.Range("M71").Value = -71256
$ws_CRP.Range("H122").Value = 1251.6666
$ws_CRP.Range("I122").Value = 1364.3
$ws_CRP.Range("J122").Value = 1026.4
$ws_CRP.Range("K122").Value = 4092.9
$ws_CRP.Range("L122").Value = 3079.2
$ws_CRP.Range("M122").Value = -1642.9
$ws_CRP.Range("N122").Value = -7979.200000000001
$ws_CRP.Range("H136").Value = 1819.48
$ws_CRP.Range("I136").Value = 1702.9546
$ws_CRP.Range("K136").Value = 5108.8638
$ws_CRP.Range("M136").Value = -2558.8638
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 1627.52
$ws_CUL.Range("I5").Value = 2409
$ws_CUL.Range("J5").Value = 906.1539
$ws_CUL.Range("K5").Value = 7227
$ws_CUL.Range("L5").Value = 2718.4617
$ws_CUL.Range("M5").Value = -7115
$ws_CUL.Range("N5").Value = -2942.4617
$ws_CUL.Range("H107").Value = 591.4375
$ws_CUL.Range("J107").Value = 392.5
$ws_CUL.Range("L107").Value = 1177.5
$ws_CUL.Range("N107").Value = -5017.5
$ws_CUL.Range("H122").Value = 758.8
$ws_CUL.Range("I122").Value = 498.27777
$ws_CUL.Range("J122").Value = 1428.7142
$ws_CUL.Range("K122").Value = 4484.49993
$ws_CUL.Range("L122").Value = 12858.4278
$ws_CUL.Range("M122").Value = -2034.49993
$ws_CUL.Range("N122").Value = -17758.4278
$ws_CUL.Range("H123").Value = 1100
$ws_CUL.Range("J123").Value = 0
$ws_CUL.Range("L123").Value = 0
$ws_CUL.Range("N123").ClearContents()
$ws_CUL.Range("H125").Value = 3009.4119
$ws_CUL.Range("I125").Value = 1944
$ws_CUL.Range("J125").Value = 3453.3333
$ws_CUL.Range("K125").Value = 5832
$ws_CUL.Range("L125").Value = 10359.9999
$ws_CUL.Range("M125").Value = -912
$ws_CUL.Range("N125").Value = -20199.9999
$ws_CUL.Range("H135").Value = 1627.52
$ws_CUL.Range("I135").Value = 2409
$ws_CUL.Range("J135").Value = 906.1539
$ws_CUL.Range("K135").Value = 21681
$ws_CUL.Range("L135").Value = 8155.3851
$ws_CUL.Range("M135").Value = -19146
$ws_CUL.Range("N135").Value = -13225.3851
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H64").Value = 20000
$ws_GSM.Range("I64").Value = 20000
$ws_GSM.Range("K64").Value = 20000
$ws_GSM.Range("M64").Value = -19752
$ws_GSM.Range("H67").Value = 20000
$ws_GSM.Range("I67").Value = 20000
$ws_GSM.Range("K67").Value = 20000
$ws_GSM.Range("M67").Value = -19142
$ws_GSM.Range("H69").Value = 147444.44
$ws_GSM.Range("J69").Value = 147444.44
$ws_GSM.Range("L69").Value = 147444.44
$ws_GSM.Range("N69").Value = -148942.44
$ws_GSM.Range("H72").Value = 147444.44
$ws_GSM.Range("J72").Value = 147444.44
$ws_GSM.Range("L72").Value = 442333.32
$ws_GSM.Range("N72").Value = -449821.32
$ws_GSM.Range("H102").Value = 1814.2424
$ws_GSM.Range("I102").Value = 1847.16
$ws_GSM.Range("K102").Value = 1847.16
$ws_GSM.Range("M102").Value = -225.1600000000001
$ws_GSM.Range("H109").Value = 10277
$ws_GSM.Range("J109").Value = 10277
$ws_GSM.Range("L109").Value = 10277
$ws_GSM.Range("N109").Value = -12357
$ws_GSM.Range("H122").Value = 3231.9167
$ws_GSM.Range("I122").Value = 3358.8
$ws_GSM.Range("J122").Value = 2597.5
$ws_GSM.Range("K122").Value = 10076.4
$ws_GSM.Range("L122").Value = 7792.5
$ws_GSM.Range("M122").Value = -7626.400000000001
$ws_GSM.Range("N122").Value = -12692.5
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 6454
$ws_LTW.Range("I40").Value = 7735.6
$ws_LTW.Range("K40").Value = 7735.6
$ws_LTW.Range("M40").Value = -7599.6
$ws_LTW.Range("H122").Value = 22733028
$ws_LTW.Range("I122").Value = 27783700
$ws_LTW.Range("J122").Value = 5000
$ws_LTW.Range("K122").Value = 83351100
$ws_LTW.Range("L122").Value = 15000
$ws_LTW.Range("M122").Value = -83348650
$ws_LTW.Range("N122").Value = -19900
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H64").Value = 0
$ws_WVR.Range("J64").Value = 0
$ws_WVR.Range("L64").Value = 0
$ws_WVR.Range("N64").ClearContents()
$ws_WVR.Range("H67").Value = 0
$ws_WVR.Range("J67").Value = 0
$ws_WVR.Range("L67").Value = 0
$ws_WVR.Range("N67").ClearContents()
$ws_WVR.Range("H70").Value = 75833.336
$ws_WVR.Range("I70").Value = 136666.67
$ws_WVR.Range("K70").Value = 136666.67
$ws_WVR.Range("M70").Value = -136351.67
$ws_WVR.Range("H73").Value = 75833.336
$ws_WVR.Range("I73").Value = 136666.67
$ws_WVR.Range("K73").Value = 136666.67
$ws_WVR.Range("M73").Value = -135574.67
$ws_WVR.Range("H122").Value = 83343336
$ws_WVR.Range("I122").Value = 125005000
$ws_WVR.Range("J122").Value = 20000
$ws_WVR.Range("K122").Value = 375015000
$ws_WVR.Range("L122").Value = 60000
$ws_WVR.Range("M122").Value = -375012550
$ws_WVR.Range("N122").Value = -64900
